$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.053638458251953
$ws.Range("B1").Value = 1.574147582054138
$ws.Range("C1").Value = 2.800167322158813
$ws.Range("D1").Value = 1.534817695617676
$ws.Range("E1").Value = 0.8238430023193359
